# Apply "group by/having/order by as separate nodes" edit to Project Notes.xlsx
#
# Touches two sheets:
#  - "Meeting Notes" (sheet1) -> no actual content change (string table
#    dedup shift only in the source diff; text stays identical)
#  - "SQL Parser" (sheet2)    -> restructure the small notes grid: drop the
#    "case" placeholder cells, add a bunch of new annotation/question notes,
#    and rework the "TPCH Correlated" test-matrix block (rows 18-25)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL Parser")

# ---------------------------------------------------------------------
# Rows 1-14: small "RA operator questions" grid
# ---------------------------------------------------------------------

# The two "case" placeholder cells (highlighted, style s="13") lose their
# text but keep their formatting.
$ws.Range("E3").ClearContents()
$ws.Range("D4").ClearContents()

# F4 gains the highlighted fill used elsewhere in column C (style s="10");
# grab it via copy/paste-special so we reuse the existing style record.
$ws.Range("C9").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("F4").Value = "like"

# New notes added in column F / I for rows 8 and 9
$ws.Range("F8").Value = "schema"
$ws.Range("I8").Value = "What makes Q2 harder than Q1? Same RA tree"
$ws.Range("F9").Value = "case expression"

# Row 9's leftover "case" label cell (A9, no special formatting) disappears
$ws.Range("A9").ClearContents()

# New "Todos" list in column I, rows 11-15, matching the fill already used
# in column C for rows 9-14 / row 1's header style where applicable.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("I11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats (style 8)
$ws.Application.CutCopyMode = $false
$ws.Range("I11").Value = "Todos"

$ws.Range("A2").Copy() | Out-Null
$ws.Range("I12:I14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats (style 9)
$ws.Application.CutCopyMode = $false
$ws.Range("I12").Value = "group by as own operator"
$ws.Range("I13").Value = "order by = sort operator (unter projection)"
$ws.Range("I14").Value = "having operator (like where)"

# New row 15
$ws.Range("D15").Value = " "
$ws.Range("I15").Value = "add step where scan whole query and map attributes to relations (for tpch)"

# ---------------------------------------------------------------------
# Rows 18-25: "TPCH Correlated" test matrix block
# ---------------------------------------------------------------------

# Shift the existing "Neumann Q1, Q2" / "Q1" / "Q2" column from C to D,
# which also clears out column C completely (no leftover style/cell).
$ws.Range("C18:C20").Cut($ws.Range("D18:D20")) | Out-Null
$ws.Range("C18:C20").Clear() | Out-Null

# New column B header + values, matching the bold/filled header style (12)
# used by A18/D18 for the header row, and plain style for the data rows.
$ws.Range("A18").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats (style 12)
$ws.Application.CutCopyMode = $false
$ws.Range("B18").Value = "missing"

$ws.Range("B19").Value = "like"
$ws.Range("B20").Value = "date,interval,exists"
$ws.Range("B22").Value = "in, double nested,date,interval"
$ws.Range("B23").Value = "exists/not exists"
$ws.Range("B24").Value = "in list,substring,not exists"

# Column A (numeric query-ids) rows 19-24 now use left-aligned style (6)
$ws.Range("A19:A24").HorizontalAlignment = -4131  # xlLeft

# New row 25
$ws.Range("A25").Value = "all"
$ws.Range("C9").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null  # xlPasteFormats (style 10)
$ws.Application.CutCopyMode = $false
$ws.Range("B25").Value = "detect correlated subquery without alias"

# ---------------------------------------------------------------------
# Column widths / view state
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.7109375
$ws.Application.ActiveWindow.Zoom = 110
$ws.Range("E11").Select()
